$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MPXV_Metadata")

# Remove the now-unused column O (15) — shifts everything from P onward left by one.
$ws.Columns.Item(15).Delete()

# Add the missing bioproject ID columns (test_group / PRJNA849962) for each data row.
$ws.Range("C3").Value = "test_group"
$ws.Range("D3").Value = "PRJNA849962"
$ws.Range("C4").Value = "test_group"
$ws.Range("D4").Value = "PRJNA849962"
$ws.Range("C5").Value = "test_group"
$ws.Range("D5").Value = "PRJNA849962"
$ws.Range("C6").Value = "test_group"
$ws.Range("D6").Value = "PRJNA849962"
$ws.Range("C7").Value = "test_group"
$ws.Range("D7").Value = "PRJNA849962"

# Update the saved view state to match: no frozen top-left cell, selection at F7.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F7").Select()
